# Add a new student record (Bob) as row 9 of Sheet1.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 9

# Cells whose text looks numeric/date-like need to be forced to Text format
# *before* assignment, otherwise Excel auto-converts them to a real number /
# date serial. We reset the style back to Normal afterwards so the cells end
# up with no explicit style (matching the rest of the sheet) while still
# being stored as text.
$textForcedCells = @("A$row", "C$row", "G$row", "H$row", "K$row")
foreach ($addr in $textForcedCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("A$row").Value = "1234567890123"
$ws.Range("B$row").Value = "Bob"
$ws.Range("C$row").Value = "1234509876"
$ws.Range("D$row").Value = "A"
$ws.Range("E$row").Value = "Yes"
$ws.Range("F$row").Value = "3 months"
$ws.Range("G$row").Value = "2025-01-11"
$ws.Range("H$row").Value = "2025-04-23"
$ws.Range("I$row").Value = "Boom Company"
$ws.Range("J$row").Value = "Campus"
$ws.Range("K$row").Value = "123000"
$ws.Range("L$row").Value = "Industry"
$ws.Range("M$row").Value = "Chennai"
$ws.Range("N$row").Value = "No"
$ws.Range("O$row").Value = ""
$ws.Range("P$row").Value = ""
$ws.Range("Q$row").Value = ""
$ws.Range("R$row").Value = "Yes"

foreach ($addr in $textForcedCells) {
    $ws.Range($addr).Style = "Normal"
}
